$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty G11 / H11 cells with value 5
# (their style, s=2, is already correct and unchanged)
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 5

# Update F19 from 2 to 5, and switch its look from the "green" style (s=4)
# to the plain style (s=2) used by the rest of the already-answered cells
# in that row (copy E19's format onto F19, then set the new value).
$ws.Range("E19").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value = 5

# Update the active selection to match the newly edited cells
$ws.Range("G11:H11").Select()
